{"js": "// Replace the placeholder paragraph \">>>  your stuff after this line >>>\"\n// with \">>>  Version control for assignment 1 part 2, Qiang Wang 23902684>>>\"\n// using three clean runs (matching the target OOXML) and no leftover\n// proofErr markers.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Locate the paragraph that still contains the placeholder text.\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text.indexOf(\">>>\") !== -1 && text.indexOf(\"your stuff after this line\") !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (target) {\n  const ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">&gt;&gt;&gt;  </w:t></w:r>\n            <w:r><w:t>Version control for assignment 1 part 2, Qiang Wang 23902684</w:t></w:r>\n            <w:r><w:t>&gt;&gt;&gt;</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\n  const range = target.getRange(\"Whole\");\n  range.insertOoxml(ooxml, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace the placeholder paragraph \">>>  your stuff after this line >>>\"\n# with \">>>  Version control for assignment 1 part 2, Qiang Wang 23902684>>>\"\n# using three clean runs (matching the target OOXML) and no leftover\n# proofErr markers.\n\n$d = $word.ActiveDocument\n\n$target = $null\nforeach ($para in $d.Paragraphs) {\n    $t = $para.Range.Text\n    if ($t -like \"*>>>*\" -and $t -like \"*your stuff after this line*\") {\n        $target = $para\n        break\n    }\n}\n\nif ($target -ne $null) {\n    $ooxml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r><w:t xml:space=\"preserve\">&gt;&gt;&gt;  </w:t></w:r>\n            <w:r><w:t>Version control for assignment 1 part 2, Qiang Wang 23902684</w:t></w:r>\n            <w:r><w:t>&gt;&gt;&gt;</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n    $target.Range.InsertXML($ooxml)\n}\n"}
